# Updated Global_M2 for easier usage.
# Apply updated Hungary M2 / FX / computed-value figures for the existing
# monthly rows, and append two new rows (388, 389) with the latest data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows with revised source data --------------------
$ws.Range("B358").Value = 30734900000000
$ws.Range("D358").Value = 99090818231.3513

$ws.Range("B359").Value = 31322400000000
$ws.Range("D359").Value = 99561353320.36465

$ws.Range("B361").Value = 33495800000000
$ws.Range("D361").Value = 112815417570.4259

$ws.Range("B362").Value = 32988400000000
$ws.Range("D362").Value = 111980718965.3417

$ws.Range("B363").Value = 33552400000000
$ws.Range("D363").Value = 111796614687.4583

$ws.Range("B364").Value = 33970200000000
$ws.Range("D364").Value = 110071285075.4974

$ws.Range("B365").Value = 33881500000000
$ws.Range("D365").Value = 113126878130.217

$ws.Range("B366").Value = 33795100000000
$ws.Range("D366").Value = 118994736008.1689

$ws.Range("B367").Value = 34176400000000
$ws.Range("D367").Value = 115377006566.1766

$ws.Range("B368").Value = 34635200000000
$ws.Range("D368").Value = 114646232270.2372

$ws.Range("B369").Value = 34976400000000
$ws.Range("D369").Value = 118284995400.6818

$ws.Range("B370").Value = 35509700000000
$ws.Range("D370").Value = 114397595407.3053

$ws.Range("B371").Value = 36350700000000
$ws.Range("D371").Value = 116727517942.2956

$ws.Range("B372").Value = 37653400000000
$ws.Range("D372").Value = 117063267526.8149

$ws.Range("B373").Value = 38869800000000
$ws.Range("D373").Value = 119742338545.7099

$ws.Range("B375").Value = 39527800000000
$ws.Range("D375").Value = 119368847013.3478

$ws.Range("B376").Value = 39600300000000
$ws.Range("D376").Value = 119319342543.3509

$ws.Range("B377").Value = 39874100000000
$ws.Range("D377").Value = 111147317073.1707

$ws.Range("B380").Value = 41209500000000
$ws.Range("D380").Value = 104230214735.5642

$ws.Range("B381").Value = 42040500000000
$ws.Range("D381").Value = 105465104610.9076

$ws.Range("B382").Value = 41918300000000
$ws.Range("D382").Value = 97100532777.39171

$ws.Range("B387").Value = 40513200000000
$ws.Range("D387").Value = 113226720625.1432

# --- Append two new monthly rows (388, 389) ----------------------------
# Copy the formatting of the last existing row down onto the new rows so
# the new date cells (column A) pick up the same style (bold, centered,
# bordered, custom date number format) used throughout the table.
$ws.Range("A387:D387").Copy()
$ws.Range("A388:D389").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A388").Value = 44986
$ws.Range("B388").Value = 40410200000000
$ws.Range("C388").Value = 0.00285257873117298
$ws.Range("D388").Value = 115273277042.4464

$ws.Range("A389").Value = 45017
$ws.Range("B389").Value = 39781400000000
$ws.Range("C389").Value = 0.002952770436862386
$ws.Range("D389").Value = 117465341856.9973
